# ptxdiffconc.xlsx: switch the "latency" column (L) from a computed
# 1/K value to hard numbers taken from the median (not mean) of all
# traces, add a "stdev" column (S) and a "% difference" column header
# for R, plus summary cells R18/Q19. See commit message:
# "using the median of all traces rather than just the mean trace
#  leads to robustness in the model"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates / additions ---------------------------------
# Shared-string table order matters for byte-identical output: the
# original file appends new unique strings in first-use order, so we
# add "stdev" (S1), then "% difference" (R1), then "latency" (L1).
$ws.Range("S1").Value2 = "stdev"
$ws.Range("R1").Value2 = "% difference"
$ws.Range("L1").Value2 = "latency"

# --- New hard-coded "latency" values for column L (rows 2-17) -------
# (replaces the old `=1/K#` formula with literal numbers)
$L = @{
    2  = 250
    3  = 128.5
    4  = 281
    5  = 139
    6  = 333
    7  = 83
    8  = 300.5
    9  = 83
    10 = 280
    11 = 141
    12 = 196
    13 = 82
    14 = 297
    15 = 94
    16 = 275
    17 = 79
}
foreach ($r in $L.Keys) {
    $ws.Range("L$r").Value2 = $L[$r]
}

# --- New "stdev" values for column S (rows 2-17) ---------------------
$S = @{
    2  = 56.9748
    3  = 22.8279
    4  = 13.6158
    5  = 105.6248
    6  = 161.639
    7  = 160.2904
    8  = 222.4139
    9  = 67.8848
    10 = 79.4504
    11 = 116.1221
    12 = 133.9883
    13 = 101.732
    14 = 86.0576
    15 = 105.4444
    16 = 153.5434
    17 = 83.5901
}
foreach ($r in $S.Keys) {
    $ws.Range("S$r").Value2 = $S[$r]
}

# --- New summary formulas --------------------------------------------
$ws.Range("R18").Formula = "=AVERAGE(R2:R17)"
$ws.Range("Q19").Formula = "=(Q18-P18)*100/Q18"

# --- Bold the outlier "% difference" cells (R5, R6, R11) -------------
$ws.Range("R5").Font.Bold = $true
$ws.Range("R6").Font.Bold = $true
$ws.Range("R11").Font.Bold = $true

# --- Column widths -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15.59
$ws.Columns.Item(18).ColumnWidth = 10.42

# --- View state: scroll position + selection --------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("E18").Select()
